# Edit script: update Solution Briefing presentation
#  1. Slide 1 - bump the presenter/date line to November 24, 2025
#  2. Slide 9 - Investment Summary table:
#     - Professional Services row zeroed out (removed from infra costs)
#     - "Software" label expanded to "Software Licenses"
#     - TOTAL row recalculated to drop the Professional Services amounts

$p = $ppt.ActivePresentation

# --- Slide 1: title slide presenter/date line ---------------------------
# (use TextRange2.Characters so we rewrite the existing run's <a:t> only,
#  instead of materializing a brand-new <a:rPr lang="en-US"/> on a run that
#  did not have explicit run properties before the edit)
$s1 = $p.Slides.Item(1)
$dateRange = $s1.Shapes.Item(3).TextFrame2.TextRange
$dateRange.Characters(1, $dateRange.Length).Text = "['Presenter Name'] | November 24, 2025"

# --- Slide 9: Investment Summary table -----------------------------------
$s9 = $p.Slides.Item(9)
$tbl = $s9.Shapes.Item(2).Table

# Professional Services row (row 2): Year 1 List / Year 1 Net / 3-Year Total -> $0
$tbl.Cell(2, 2).Shape.TextFrame.TextRange.Text = "$0"
$tbl.Cell(2, 4).Shape.TextFrame.TextRange.Text = "$0"
$tbl.Cell(2, 7).Shape.TextFrame.TextRange.Text = "$0"

# Software row (row 4): rename category label
$tbl.Cell(4, 1).Shape.TextFrame.TextRange.Text = "Software Licenses"

# TOTAL row (row 6): recalculated totals without Professional Services
$tbl.Cell(6, 2).Shape.TextFrame.TextRange.Text = "$157,720"
$tbl.Cell(6, 4).Shape.TextFrame.TextRange.Text = "$141,724"
$tbl.Cell(6, 7).Shape.TextFrame.TextRange.Text = "$223,964"
